# The workbook contains four lookup-table sheets (Fuel_to_Code, VehFuel_to_Code,
# Tech_to_Code, Dem_to_Code) that each carry a column header literally reading
# "Plain English". Rename that header to "Plain_English" (underscore) on every
# sheet that uses it, leaving every other cell value untouched.

$wb = $excel.ActiveWorkbook

$targets = @("Fuel_to_Code", "VehFuel_to_Code", "Tech_to_Code", "Dem_to_Code")

foreach ($name in $targets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Cells.Replace("Plain English", "Plain_English")
}
